$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: replace the first occurrence of $OldText inside a TextRange with
# $NewText while preserving the surrounding run formatting (uses
# TextRange.Characters(start,len) so only the matched span is touched).
# NOTE: positional parameters only -- named parameter binding (-Foo bar) is
# not reliable in this PowerShell host.
# ---------------------------------------------------------------------------
function Replace-InTextRange {
    param($TextRange, $OldText, $NewText)
    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -ge 0) {
        $sub = $TextRange.Characters($idx + 1, $OldText.Length)
        $sub.Text = $NewText
        return $true
    }
    return $false
}

# Update the "Date Placeholder" shape (the auto date field) on a shape
# collection (slide master or a custom layout) from the old cached date to
# the new one.
function Update-DatePlaceholder {
    param($Shapes, $OldDate, $NewDate)
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $shp = $Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            Replace-InTextRange $shp.TextFrame.TextRange $OldDate $NewDate | Out-Null
        }
    }
}

$oldDate = "2020-05-11"
$newDate = "2020-05-18"

# Refresh the cached date on the slide master ...
Update-DatePlaceholder $p.SlideMaster.Shapes $oldDate $newDate

# ... and on every slide layout that derives from it.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes $oldDate $newDate
}

# ---------------------------------------------------------------------------
# Front matter text updates on slide 1, "Content Placeholder 4" shape.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$shape = $slide.Shapes.Item("Content Placeholder 4")
$tr = $shape.TextFrame.TextRange

Replace-InTextRange $tr "Marcos" "Kyle" | Out-Null
Replace-InTextRange $tr "Use the chat button at the bottom of your screen" "Use the chat button at the bottom of your screen (try and keep questions concise)" | Out-Null
Replace-InTextRange $tr ". She will do her best to get to your question at the end of the seminar." ". He will do his best to get to your question at the end of the seminar." | Out-Null
